# Reorders rows within five table blocks on the active sheet so that each
# section's group-label row (previously the first row of the block, with
# only column A populated) moves to become the LAST row of its block, and
# the data rows that followed it each shift up by one row.
#
# Blocks affected (1-based row numbers, inclusive): 7-8, 12-14, 17-20,
# 22-26, 27-32. Columns A:L.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 200

$blocks = @(
    @{ Start = 7;  End = 8  },
    @{ Start = 12; End = 14 },
    @{ Start = 17; End = 20 },
    @{ Start = 22; End = 26 },
    @{ Start = 27; End = 32 }
)

foreach ($block in $blocks) {
    $start = $block.Start
    $end = $block.End
    $n = $end - $start + 1

    $scratchStart = $scratchRow
    $scratchEnd = $scratchRow + $n - 1

    # 1) Stash the whole block (in original order) into the scratch area.
    $ws.Range("A$start`:L$end").Copy()
    $ws.Range("A$scratchStart").PasteSpecial()

    # 2) Move the data rows (everything after the group-label row) up by one,
    #    taking them from the scratch copy.
    $dataScratchStart = $scratchStart + 1
    $ws.Range("A$dataScratchStart`:L$scratchEnd").Copy()
    $ws.Range("A$start").PasteSpecial()

    # 3) The last row of the block becomes the group-label row, so clear any
    #    leftover data in columns B:L before writing it.
    $ws.Range("B$end`:L$end").ClearContents()
    $ws.Range("A$scratchStart`:L$scratchStart").Copy()
    $ws.Range("A$end").PasteSpecial()

    # 4) Clean up the scratch area.
    $ws.Range("A$scratchStart`:L$scratchEnd").Clear()
}

$excel.CutCopyMode = $false
